# Applies the "Finished Letter to Mr. Judan" edit to the Questionnaire document.
#
# Strategy: the runs we need to touch sit next to sibling runs that share
# identical rPr (e.g. a "[]" checkbox run followed by a label run). Rewriting
# text through Find/Replace or Range.Text normalizes/merges such runs, so
# every textual change below is done with Range.InsertXML against a
# precisely bounded Range (computed from the pristine document) which
# leaves neighboring runs completely untouched. Edits are applied from the
# end of the document towards the start so earlier offsets stay valid.
#
# NOTE: this PowerShell host does not support nested function-call
# expressions as call arguments (e.g. `Foo 1 (Bar)`), so every helper
# result is assigned to a variable before being passed along.

$d = $word.ActiveDocument
$apos = [char]0x2019
$box = [char]0x25A1

function New-RunXml([string]$text, [bool]$preserve) {
    $space = ""
    if ($preserve) { $space = ' xml:space="preserve"' }
    $xml = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t' + $space + '>' + $text + '</w:t></w:r>'
    return $xml
}

function Wrap-Pkg([string]$bodyXml) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $xml
}

function Set-RangeXml([int]$start, [int]$end, [string]$innerXml) {
    $body = '<w:body><w:p>' + $innerXml + '</w:p></w:body>'
    $pkg = Wrap-Pkg $body
    $rng = $d.Range($start, $end)
    $rng.InsertXML($pkg)
}

# --- 1. "Peers" paragraph (15): change text, then add a brand-new empty
#        paragraph right after it (before the doc's existing trailing
#        empty paragraph). Handled back-to-front: insert the new empty
#        paragraph first (at the end of the Peers paragraph), then fix the
#        "Peers" run text itself.
$peersPara = $d.Paragraphs.Item(15)
$peersEnd = $peersPara.Range.End
$newEmptyParaBody = '<w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p></w:body>'
$newEmptyParaXml = Wrap-Pkg $newEmptyParaBody
$insAt = $d.Range($peersEnd, $peersEnd)
$insAt.InsertXML($newEmptyParaXml)

$peersRunXml = New-RunXml " My Peers are in the college of my choice" $true
Set-RangeXml 445 450 $peersRunXml

# --- 2. "[]Tuition Fees" (14): split the single run into "[]" + the new
#        sentence.
$boxRunXml = New-RunXml $box $false
$tuitionRunXml = New-RunXml " The Tuition Fee is higher here compared to the college of my choice" $true
$tuitionXml = $boxRunXml + $tuitionRunXml
Set-RangeXml 419 432 $tuitionXml

# --- 3. "The Location of the School (Commute, Bus-service)" (13)
$locationXml = New-RunXml " The Availability of the school is a hassle " $true
Set-RangeXml 358 407 $locationXml

# --- 4. "Quality of Training" (12)
$qualityXml = New-RunXml " The Quality of Training does not suit to my needs" $true
Set-RangeXml 326 345 $qualityXml

# --- 5. "Reputation of the School" (11)
$reputationXml = New-RunXml " The Reputation of the School might be a hindrance for getting a job" $true
Set-RangeXml 289 313 $reputationXml

# --- 6. "Program Offered & Specialization" (10)
$programText = "The Program of my choice isn" + $apos + "t offered here"
$programXml = New-RunXml $programText $false
Set-RangeXml 244 276 $programXml

# --- 7. "If No: Check the reasons why you do not pursue education here" (9)
#        becomes 5 runs; leading whitespace is dropped entirely.
$run1 = New-RunXml "If No: Check the reasons" $false
$run2 = New-RunXml " that apply to you" $true
$run3 = New-RunXml " why you do not pursue " $true
$run4 = New-RunXml "further " $true
$run5 = New-RunXml "education here" $false
$para9Xml = $run1 + $run2 + $run3 + $run4 + $run5
Set-RangeXml 159 231 $para9Xml

# --- 8. Remove the old (mid-document) "_GoBack" bookmark; it is recreated
#        at the top of the document in step 9.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- 9. Insert the new leading "Directions: ..." paragraph (with the
#        relocated "_GoBack" bookmark) before the first paragraph.
$introP = '<w:p><w:pPr><w:spacing w:before="240" w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Directions: For the following items, please fill-up the items below by checking your appropriate responses. All your answers will be </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>treated with utmost confidence.</w:t></w:r></w:p>'
$introBody = '<w:body>' + $introP + '</w:body>'
$introPkg = Wrap-Pkg $introBody
$startRng = $d.Range(0, 0)
$startRng.InsertXML($introPkg)

Write-Output "done"
